# Actualización automática 2025-09-11 09:20:09
#
# Updates sales figures for CHASI PASTO ANGEL NOLBERTO (row 17) on the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, and rolls those changes
# through the dependent total / count rows on "VENTAS POR GRUPO" (row 60),
# "VENTA MENSUAL" (row 60) and the "CUMPLIMIENTO MENSUAL" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Cells.Item(17, 4).Value2 = 950.4              # D17 240X80 PORCELANATO
$wsGrupo.Cells.Item(17, 8).Value2 = 1188                # H17 INODOROS
$wsGrupo.Cells.Item(17, 9).Value2 = 559.6799999999999   # I17 LAVABOS
$wsGrupo.Cells.Item(17, 13).Value2 = 3758.1             # M17 PORCELANATO

# Row 60 holds "<n> de 58" counts of non-zero entries per column; the four
# columns touched above each gain one more non-zero row.
$wsGrupo.Cells.Item(60, 4).Value2 = "1 de 58"   # D60
$wsGrupo.Cells.Item(60, 8).Value2 = "1 de 58"   # H60
$wsGrupo.Cells.Item(60, 9).Value2 = "1 de 58"   # I60
$wsGrupo.Cells.Item(60, 13).Value2 = "3 de 58"  # M60

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Cells.Item(17, 6).Value2 = 6456.18   # F17 septiembre
$wsMensual.Cells.Item(60, 6).Value2 = 12321.45  # F60 total septiembre

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3  - 240X80 PORCELANATO
$wsCumplimiento.Cells.Item(3, 4).Value2 = 950.4
$wsCumplimiento.Cells.Item(3, 5).Value2 = 4554.21890386263
$wsCumplimiento.Cells.Item(3, 6).Value2 = 0.1726550042062126

# Row 6  - INODOROS
$wsCumplimiento.Cells.Item(6, 4).Value2 = 1188
$wsCumplimiento.Cells.Item(6, 5).Value2 = -338.15580317996
$wsCumplimiento.Cells.Item(6, 6).Value2 = 1.397903291503639

# Row 7  - LAVABOS
$wsCumplimiento.Cells.Item(7, 4).Value2 = 559.6799999999999
$wsCumplimiento.Cells.Item(7, 5).Value2 = 149.688813030059
$wsCumplimiento.Cells.Item(7, 6).Value2 = 0.7889830927431595

# Row 12 - PORCELANATO
$wsCumplimiento.Cells.Item(12, 4).Value2 = 4056.41
$wsCumplimiento.Cells.Item(12, 5).Value2 = 28348.39
$wsCumplimiento.Cells.Item(12, 6).Value2 = 0.1251792944255172

# Row 15 - TOTAL
$wsCumplimiento.Cells.Item(15, 4).Value2 = 12321.45
$wsCumplimiento.Cells.Item(15, 5).Value2 = 38162.31705102521
$wsCumplimiento.Cells.Item(15, 6).Value2 = 0.2440675630950124

# Column F (CUMPLIMIENTO) narrows from 26 to 24 characters wide.
# ColumnWidth is expressed in characters and is offset from the raw
# stored OOXML width by 5/6, so ask for (24 - 5/6) to land on 24 exactly.
$wsCumplimiento.Columns.Item(6).ColumnWidth = 23.166666666666668
